# Map the Mongo _id onto row 2 for Excel visualization, and update the
# company's details (name, impact level, foundation year, description)
# plus bump the updatedAt timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "67bf38f0a9c0f101fb8c8295"
$ws.Range("B2").Value = "Pollo Grangero"
$ws.Range("C2").Value = "Medio"
$ws.Range("D2").Value = 2010
$ws.Range("F2").Value = "Empresa de comida."
$ws.Range("H2").Value = 45714.430585925926
